$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Angpt1"
$ws.Cells.Item(2,3).Value = "Itgb1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.110507
$ws.Cells.Item(2,8).Value = 0.331521
$ws.Cells.Item(2,9).Value = 0.004605687348208628
$ws.Cells.Item(2,10).Value = 0.004605687348208628
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 153.5290173333333
$ws.Cells.Item(2,14).Value = 460.587052
$ws.Cells.Item(2,15).Value = 0.3172206968818489
$ws.Cells.Item(2,16).Value = 0.317220696881849
$ws.Cells.Item(2,17).Value = 16.96603111845467
$ws.Cells.Item(2,18).Value = 152.694280066092
$ws.Cells.Item(2,19).Value = 0.001461019350218656
$ws.Cells.Item(2,20).Value = 0.001461019350218656

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Angpt1"
$ws.Cells.Item(3,3).Value = "Itgb1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.110507
$ws.Cells.Item(3,8).Value = 0.331521
$ws.Cells.Item(3,9).Value = 0.004605687348208628
$ws.Cells.Item(3,10).Value = 0.004605687348208628
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 168.7997026666667
$ws.Cells.Item(3,14).Value = 506.3991080000001
$ws.Cells.Item(3,15).Value = 0.3487728915577651
$ws.Cells.Item(3,16).Value = 0.3487728915577651
$ws.Cells.Item(3,17).Value = 18.65354874258534
$ws.Cells.Item(3,18).Value = 167.881938683268
$ws.Cells.Item(3,19).Value = 0.001606338894045738
$ws.Cells.Item(3,20).Value = 0.001606338894045738

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Angpt1"
$ws.Cells.Item(4,3).Value = "Itgb1"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.110507
$ws.Cells.Item(4,8).Value = 0.331521
$ws.Cells.Item(4,9).Value = 0.004605687348208628
$ws.Cells.Item(4,10).Value = 0.004605687348208628
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 68.09032333333333
$ws.Cells.Item(4,14).Value = 204.27097
$ws.Cells.Item(4,15).Value = 0.1406878008722904
$ws.Cells.Item(4,16).Value = 0.1406878008722904
$ws.Cells.Item(4,17).Value = 7.524457360596667
$ws.Cells.Item(4,18).Value = 67.72011624537001
$ws.Cells.Item(4,19).Value = 0.0006479640245248025
$ws.Cells.Item(4,20).Value = 0.0006479640245248027

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Angpt1"
$ws.Cells.Item(5,3).Value = "Itgb1"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.110507
$ws.Cells.Item(5,8).Value = 0.331521
$ws.Cells.Item(5,9).Value = 0.004605687348208628
$ws.Cells.Item(5,10).Value = 0.004605687348208628
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 93.562673
$ws.Cells.Item(5,14).Value = 280.688019
$ws.Cells.Item(5,15).Value = 0.1933186106880956
$ws.Cells.Item(5,16).Value = 0.1933186106880956
$ws.Cells.Item(5,17).Value = 10.339330305211
$ws.Cells.Item(5,18).Value = 93.053972746899
$ws.Cells.Item(5,19).Value = 0.0008903650794194313
$ws.Cells.Item(5,20).Value = 0.0008903650794194313

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Angpt1"
$ws.Cells.Item(6,3).Value = "Itgb1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 22.14783133333333
$ws.Cells.Item(6,8).Value = 66.443494
$ws.Cells.Item(6,9).Value = 0.9230726249214253
$ws.Cells.Item(6,10).Value = 0.9230726249214253
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 153.5290173333333
$ws.Cells.Item(6,14).Value = 460.587052
$ws.Cells.Item(6,15).Value = 0.3172206968818489
$ws.Cells.Item(6,16).Value = 0.317220696881849
$ws.Cells.Item(6,17).Value = 3400.334780671076
$ws.Cells.Item(6,18).Value = 30603.01302603969
$ws.Cells.Item(6,19).Value = 0.2928177413501321
$ws.Cells.Item(6,20).Value = 0.2928177413501321

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Angpt1"
$ws.Cells.Item(7,3).Value = "Itgb1"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 22.14783133333333
$ws.Cells.Item(7,8).Value = 66.443494
$ws.Cells.Item(7,9).Value = 0.9230726249214253
$ws.Cells.Item(7,10).Value = 0.9230726249214253
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 168.7997026666667
$ws.Cells.Item(7,14).Value = 506.3991080000001
$ws.Cells.Item(7,15).Value = 0.3487728915577651
$ws.Cells.Item(7,16).Value = 0.3487728915577651
$ws.Cells.Item(7,17).Value = 3738.547343778151
$ws.Cells.Item(7,18).Value = 33646.92609400336
$ws.Cells.Item(7,19).Value = 0.3219427085116618
$ws.Cells.Item(7,20).Value = 0.3219427085116618

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Angpt1"
$ws.Cells.Item(8,3).Value = "Itgb1"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 22.14783133333333
$ws.Cells.Item(8,8).Value = 66.443494
$ws.Cells.Item(8,9).Value = 0.9230726249214253
$ws.Cells.Item(8,10).Value = 0.9230726249214253
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 68.09032333333333
$ws.Cells.Item(8,14).Value = 204.27097
$ws.Cells.Item(8,15).Value = 0.1406878008722904
$ws.Cells.Item(8,16).Value = 0.1406878008722904
$ws.Cells.Item(8,17).Value = 1508.052996618798
$ws.Cells.Item(8,18).Value = 13572.47696956918
$ws.Cells.Item(8,19).Value = 0.1298650576456079
$ws.Cells.Item(8,20).Value = 0.1298650576456079

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Angpt1"
$ws.Cells.Item(9,3).Value = "Itgb1"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 22.14783133333333
$ws.Cells.Item(9,8).Value = 66.443494
$ws.Cells.Item(9,9).Value = 0.9230726249214253
$ws.Cells.Item(9,10).Value = 0.9230726249214253
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 93.562673
$ws.Cells.Item(9,14).Value = 280.688019
$ws.Cells.Item(9,15).Value = 0.1933186106880956
$ws.Cells.Item(9,16).Value = 0.1933186106880956
$ws.Cells.Item(9,17).Value = 2072.210300699821
$ws.Cells.Item(9,18).Value = 18649.89270629839
$ws.Cells.Item(9,19).Value = 0.1784471174140236
$ws.Cells.Item(9,20).Value = 0.1784471174140236

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Angpt1"
$ws.Cells.Item(10,3).Value = "Itgb1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.735257333333333
$ws.Cells.Item(10,8).Value = 5.205772
$ws.Cells.Item(10,9).Value = 0.07232168773036617
$ws.Cells.Item(10,10).Value = 0.07232168773036617
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 153.5290173333333
$ws.Cells.Item(10,14).Value = 460.587052
$ws.Cells.Item(10,15).Value = 0.3172206968818489
$ws.Cells.Item(10,16).Value = 0.317220696881849
$ws.Cells.Item(10,17).Value = 266.4123532071271
$ws.Cells.Item(10,18).Value = 2397.711178864144
$ws.Cells.Item(10,19).Value = 0.02294193618149822
$ws.Cells.Item(10,20).Value = 0.02294193618149822

# Row 11
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Angpt1"
$ws.Cells.Item(11,3).Value = "Itgb1"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 1.735257333333333
$ws.Cells.Item(11,8).Value = 5.205772
$ws.Cells.Item(11,9).Value = 0.07232168773036617
$ws.Cells.Item(11,10).Value = 0.07232168773036617
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 168.7997026666667
$ws.Cells.Item(11,14).Value = 506.3991080000001
$ws.Cells.Item(11,15).Value = 0.3487728915577651
$ws.Cells.Item(11,16).Value = 0.3487728915577651
$ws.Cells.Item(11,17).Value = 292.9109219168196
$ws.Cells.Item(11,18).Value = 2636.198297251376
$ws.Cells.Item(11,19).Value = 0.02522384415205755
$ws.Cells.Item(11,20).Value = 0.02522384415205755

# Row 12
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Angpt1"
$ws.Cells.Item(12,3).Value = "Itgb1"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 1.735257333333333
$ws.Cells.Item(12,8).Value = 5.205772
$ws.Cells.Item(12,9).Value = 0.07232168773036617
$ws.Cells.Item(12,10).Value = 0.07232168773036617
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 68.09032333333333
$ws.Cells.Item(12,14).Value = 204.27097
$ws.Cells.Item(12,15).Value = 0.1406878008722904
$ws.Cells.Item(12,16).Value = 0.1406878008722904
$ws.Cells.Item(12,17).Value = 118.1542328932044
$ws.Cells.Item(12,18).Value = 1063.38809603884
$ws.Cells.Item(12,19).Value = 0.01017477920215772
$ws.Cells.Item(12,20).Value = 0.01017477920215772

# Row 13
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Angpt1"
$ws.Cells.Item(13,3).Value = "Itgb1"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 1.735257333333333
$ws.Cells.Item(13,8).Value = 5.205772
$ws.Cells.Item(13,9).Value = 0.07232168773036617
$ws.Cells.Item(13,10).Value = 0.07232168773036617
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 93.562673
$ws.Cells.Item(13,14).Value = 280.688019
$ws.Cells.Item(13,15).Value = 0.1933186106880956
$ws.Cells.Item(13,16).Value = 0.1933186106880956
$ws.Cells.Item(13,17).Value = 162.3553144495187
$ws.Cells.Item(13,18).Value = 1461.197830045668
$ws.Cells.Item(13,19).Value = 0.01398112819465268
$ws.Cells.Item(13,20).Value = 0.01398112819465268
